$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Day 1 Opening Remarks description (row 2, column E)
$ws.Range("E2").Value = "Greetings from the FAMPS and FSN Chairs"

# Update the Day 1 Wrap-up description (row 10, column E)
$ws.Range("E10").Value = "Closing from the FAMPS and FSN Chairs; Preview of Day 2"

# Update the active selection to match the saved state
$ws.Range("E15").Select()
